$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 13:09"

# Swap province labels that moved position in the sorted table:
#  - "Asturias" now ranks above "Sevilla"
#  - "Cantabria" now ranks above "Caceres"
$ws.Range("A29").Value = "Asturias"
$ws.Range("A30").Value = "Sevilla"
$ws.Range("A32").Value = "Cantabria"
$ws.Range("A33").Value = "Caceres"

# Updated case numbers (Casos totales, Casos activos, Recuperados, Muertes)
$ws.Range("B4").Value = 64333
$ws.Range("C4").Value = 38975
$ws.Range("D4").Value = 16806
$ws.Range("E4").Value = 8552

$ws.Range("B5").Value = 51733
$ws.Range("C5").Value = 23664
$ws.Range("D5").Value = 22598
$ws.Range("E5").Value = 5471

$ws.Range("B6").Value = 17716
$ws.Range("C6").Value = 7231
$ws.Range("D6").Value = 8609
$ws.Range("E6").Value = 1876

$ws.Range("B7").Value = 16237
$ws.Range("C7").Value = 5981
$ws.Range("D7").Value = 7543
$ws.Range("E7").Value = 2713

$ws.Range("B9").Value = 12287
$ws.Range("C9").Value = 8435
$ws.Range("D9").Value = 2551
$ws.Range("E9").Value = 1301

$ws.Range("B10").Value = 9184
$ws.Range("C10").Value = 7138
$ws.Range("D10").Value = 1458
$ws.Range("E10").Value = 588

$ws.Range("B13").Value = 5274
$ws.Range("C13").Value = 2953
$ws.Range("D13").Value = 1506
$ws.Range("E13").Value = 815

$ws.Range("B15").Value = 5006
$ws.Range("C15").Value = 2827
$ws.Range("D15").Value = 1695
$ws.Range("E15").Value = 484

$ws.Range("B16").Value = 3992
$ws.Range("C16").Value = 2526
$ws.Range("D16").Value = 1124
$ws.Range("E16").Value = 342

$ws.Range("B23").Value = 2900
$ws.Range("C23").Value = 2312
$ws.Range("D23").Value = 116
$ws.Range("E23").Value = 472

$ws.Range("B29").Value = 2336
$ws.Range("C29").Value = 983
$ws.Range("D29").Value = 1061
$ws.Range("E29").Value = 292

$ws.Range("B30").Value = 2329
$ws.Range("C30").Value = 459
$ws.Range("D30").Value = 1658
$ws.Range("E30").Value = 212

$ws.Range("B31").Value = 2240
$ws.Range("C31").Value = 1318
$ws.Range("D31").Value = 774
$ws.Range("E31").Value = 148

$ws.Range("B32").Value = 2232
$ws.Range("C32").Value = 1786
$ws.Range("D32").Value = 245
$ws.Range("E32").Value = 201

$ws.Range("B33").Value = 2220
$ws.Range("C33").Value = 422
$ws.Range("D33").Value = 1482
$ws.Range("E33").Value = 316
